# Outlet -Activity testcases.xlsx - sprint 106 actual and result update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 9: add "Pass" result in column G ---
$ws.Range("G9").Value = "Pass"

# --- Row 10: update Actual output (E10) text, add Result (F10) and Pass (G10) ---
$ws.Range("E10").Value = "There is a new 'Stock count amendment' pop up screen detail page. It shows the previous qty and the Updated qty."

$run1 = $ws.Range("E10").Characters(15, 25)
$run1.Font.Bold = $true
$run1.Font.Name = "Segoe UI"
$run1.Font.Size = 11
$run1.Font.Color = 5057303

$run2 = $ws.Range("E10").Characters(40, 73)
$run2.Font.Bold = $false
$run2.Font.Name = "Segoe UI"
$run2.Font.Size = 11
$run2.Font.Color = 5057303

$ws.Range("F10").Value = "As we expected the result it displayed the Previous qty and Updated qty details"

$ws.Range("G10").Value = "Pass"

# --- Update the selection / active cell shown in the saved view ---
$ws.Activate() | Out-Null
$ws.Range("G11").Select() | Out-Null
